$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report date header (B1) moves from 28.09.2024 -> 30.09.2024
$ws.Range("B1").Value = "30.09.2024"

# Payment note (F32) moves from 29.09.2024 payment -> 01.10.2024 payment
$ws.Range("F32").Value = "01.10.2024 payment "

# Updated stock figures for the 30th report
$ws.Range("C9").Value = 477733
$ws.Range("C11").Value = 22380
$ws.Range("C12").Value = 162
$ws.Range("E20").Value = 41440
$ws.Range("E21").Value = 46791
$ws.Range("E27").Value = 28299

# Update the view/selection to match where the author left off
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("F33").Select()
